# Fruta / hortaliza, semanal
# Insert 2 new weekly price rows for "Vega Modelo de Temuco - Piña" right
# after the existing row 298, pushing the former rows 299:310 down to
# 301:312 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 299:300 (existing data shifts down to 301:312).
$ws.Rows("299:300").Insert()

# --- New row 299 ------------------------------------------------------
$ws.Range("A299").Value = 10
$ws.Range("B299").Value = 'Vega Modelo de Temuco'
$ws.Range("C299").Value = 'La Araucanía'
$ws.Range("D299").Value = 44509
$ws.Range("E299").Value = 9
$ws.Range("F299").Value = 'Fruta'
$ws.Range("G299").Value = 100108
$ws.Range("H299").Value = 'Tropicales y subtropicales'
$ws.Range("I299").Value = 100108005
$ws.Range("J299").Value = 'Piña'
$ws.Range("K299").Value = 'Caramelo'
$ws.Range("L299").Value = 'Primera'
$ws.Range("M299").Value = 55
$ws.Range("N299").Value = 21000
$ws.Range("O299").Value = 21000
$ws.Range("P299").Value = 21000
$ws.Range("Q299").Value = '$/caja 12 unidades'
$ws.Range("R299").Value = 'Ecuador'
$ws.Range("S299").Value = 1750
$ws.Range("T299").Value = 12

# --- New row 300 ------------------------------------------------------
$ws.Range("A300").Value = 10
$ws.Range("B300").Value = 'Vega Modelo de Temuco'
$ws.Range("C300").Value = 'La Araucanía'
$ws.Range("D300").Value = 44509
$ws.Range("E300").Value = 9
$ws.Range("F300").Value = 'Fruta'
$ws.Range("G300").Value = 100108
$ws.Range("H300").Value = 'Tropicales y subtropicales'
$ws.Range("I300").Value = 100108005
$ws.Range("J300").Value = 'Piña'
$ws.Range("K300").Value = 'Caramelo'
$ws.Range("L300").Value = 'Segunda'
$ws.Range("M300").Value = 45
$ws.Range("N300").Value = 20000
$ws.Range("O300").Value = 20000
$ws.Range("P300").Value = 20000
$ws.Range("Q300").Value = '$/caja 14 unidades'
$ws.Range("R300").Value = 'Ecuador'
$ws.Range("S300").Value = 1429
$ws.Range("T300").Value = 14
